# Apply cell updates from cryptos list refresh (GitHub Actions data pull)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.486.94"
$ws.Range("E2").Value = "  +0.23%  "
$ws.Range("D3").Value = "1.807.82"
$ws.Range("E3").Value = "  +0.40%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'225.05"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.99%  "
$ws.Range("E6").Value = "  +4.83%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "'38.53"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +6.51%  "
$ws.Range("D9").Value = "'0.288"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.74%  "
$ws.Range("E10").Value = "  -3.00%  "
$ws.Range("D11").Value = "'0.0982"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.79%  "
$ws.Range("D12").Value = "2.069.22"
$ws.Range("E12").Value = "  +0.45%  "
$ws.Range("D13").Value = "'11.12"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.75%  "
$ws.Range("D14").Value = "1.807.29"
$ws.Range("E14").Value = "  +0.56%  "
$ws.Range("E15").Value = "  -2.06%  "
$ws.Range("D16").Value = "34.470.92"
$ws.Range("E16").Value = "  +0.28%  "
$ws.Range("E17").Value = "  -2.57%  "
$ws.Range("D18").Value = "'68.17"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.07%  "
$ws.Range("D19").Value = "'241.33"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.31%  "
$ws.Range("D20").Value = "0.0₃0768"
$ws.Range("E20").Value = "  -2.79%  "
$ws.Range("E21").Value = "  -3.52%  "
$ws.Range("E22").Value = "  -0.07%  "
$ws.Range("E23").Value = "  -1.64%  "
$ws.Range("E24").Value = "  +1.57%  "
$ws.Range("D25").Value = "'170.55"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.22%  "
$ws.Range("D26").Value = "'7.67"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.57%  "
$ws.Range("D27").Value = "'17.48"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.07%  "
$ws.Range("D28").Value = "'0.122"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.44%  "
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("D30").Value = "'1.22"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.11%  "
$ws.Range("D31").Value = "'3.77"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.52%  "
$ws.Range("D32").Value = "'0.0514"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.50%  "
$ws.Range("D33").Value = "'3.84"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.19%  "
$ws.Range("E34").Value = "  +1.19%  "
$ws.Range("D35").Value = "'0.639"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.67%  "
$ws.Range("D36").Value = "1.306.50"
$ws.Range("E36").Value = "  -6.45%  "
$ws.Range("E37").Value = "  -0.43%  "
$ws.Range("E38").Value = "  -1.66%  "
$ws.Range("E39").Value = "  -4.77%  "
$ws.Range("D40").Value = "'82.66"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.32%  "
$ws.Range("B41").Value = "WEMIXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D41").Value = "'1.22"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.41%  "
$ws.Range("B42").Value = "HuobiToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D42").Value = "'2.44"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.83%  "
$ws.Range("E43").Value = "  -0.20%  "
$ws.Range("D44").Value = "'0.948"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.97%  "
$ws.Range("E45").Value = "  +5.68%  "
$ws.Range("D46").Value = "'0.0514"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.34%  "
$ws.Range("D47").Value = "1.971.04"
$ws.Range("E47").Value = "  +0.48%  "
$ws.Range("E48").Value = "  -3.75%  "
$ws.Range("E49").Value = "  -0.05%  "
$ws.Range("D50").Value = "'102.93"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.32%  "
$ws.Range("E51").Value = "  -5.45%  "
